$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'43.557.14"
$ws.Cells.Item(2, 5).Value = "'  +2.52%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "'2.371.59"
$ws.Cells.Item(3, 5).Value = "'  +6.46%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "'  -0.25%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'318.94"
$ws.Cells.Item(5, 5).Value = "'  +8.60%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'107.83"
$ws.Cells.Item(6, 5).Value = "'  -3.28%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.639"
$ws.Cells.Item(7, 5).Value = "'  +2.05%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "'  -0.13%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.638"
$ws.Cells.Item(9, 5).Value = "'  +4.81%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'42.50"
$ws.Cells.Item(10, 5).Value = "'  -4.43%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.0932"
$ws.Cells.Item(11, 5).Value = "'  +1.56%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'8.71"
$ws.Cells.Item(12, 5).Value = "'  -1.70%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'1.05"
$ws.Cells.Item(13, 5).Value = "'  +2.22%  "

# Row 14
$ws.Cells.Item(14, 2).Value = "'TRON"
$ws.Cells.Item(14, 3).Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(14, 4).Value = "'0.106"
$ws.Cells.Item(14, 5).Value = "'  +2.37%  "

# Row 15
$ws.Cells.Item(15, 2).Value = "'Chainlink"
$ws.Cells.Item(15, 3).Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(15, 4).Value = "'16.56"
$ws.Cells.Item(15, 5).Value = "'  +9.62%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'2.727.68"

# Row 17
$ws.Cells.Item(17, 4).Value = "'2.426.53"
$ws.Cells.Item(17, 5).Value = "'  +8.73%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'43.533.58"
$ws.Cells.Item(18, 5).Value = "'  +2.51%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "'  +2.78%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "'  -1.89%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'75.20"
$ws.Cells.Item(21, 5).Value = "'  +2.94%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'3.49"
$ws.Cells.Item(22, 5).Value = "'  +0.66%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "'  +5.89%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'261.51"
$ws.Cells.Item(24, 5).Value = "'  +13.82%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'9.24"
$ws.Cells.Item(25, 5).Value = "'  -0.25%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'12.03"
$ws.Cells.Item(26, 5).Value = "'  +3.26%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'1.00"
$ws.Cells.Item(27, 5).Value = "'  +0.04%  "

# Row 28
$ws.Cells.Item(28, 2).Value = "'InjectiveProtocol"
$ws.Cells.Item(28, 3).Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(28, 4).Value = "'38.86"
$ws.Cells.Item(28, 5).Value = "'  +1.03%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "'EthereumClassic"
$ws.Cells.Item(29, 3).Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(29, 4).Value = "'22.80"
$ws.Cells.Item(29, 5).Value = "'  +8.12%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "'Toncoin"
$ws.Cells.Item(30, 3).Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(30, 4).Value = "'2.24"
$ws.Cells.Item(30, 5).Value = "'  -0.05%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'3.21"
$ws.Cells.Item(31, 5).Value = "'  -1.80%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'174.19"
$ws.Cells.Item(32, 5).Value = "'  +0.38%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'0.0920"
$ws.Cells.Item(33, 5).Value = "'  +1.97%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'5.96"
$ws.Cells.Item(34, 5).Value = "'  +4.29%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'0.132"
$ws.Cells.Item(35, 5).Value = "'  +4.02%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'4.96"
$ws.Cells.Item(36, 5).Value = "'  -5.26%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'0.0373"
$ws.Cells.Item(37, 5).Value = "'  -1.19%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'4.08"
$ws.Cells.Item(38, 5).Value = "'  -5.94%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "'  +0.06%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'2.80"
$ws.Cells.Item(40, 5).Value = "'  +15.95%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'1.49"
$ws.Cells.Item(41, 5).Value = "'  +12.45%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'71.85"
$ws.Cells.Item(42, 5).Value = "'  -0.88%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "'0.233"
$ws.Cells.Item(43, 5).Value = "'  -1.14%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "'  +0.02%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "'  -1.36%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'5.61"
$ws.Cells.Item(46, 5).Value = "'  +2.81%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "'Aave"
$ws.Cells.Item(47, 3).Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(47, 4).Value = "'111.94"
$ws.Cells.Item(47, 5).Value = "'  +8.04%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "'FraxShare"
$ws.Cells.Item(48, 3).Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(48, 4).Value = "'9.30"
$ws.Cells.Item(48, 5).Value = "'  +8.25%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "'  -1.54%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "'  +2.58%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'0.473"
$ws.Cells.Item(51, 5).Value = "'  +7.15%  "

Write-Host "Applied 103 cell updates"
